$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.030.12"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.824.07"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4664"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3661"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07236"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8607"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.87"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07555"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.767.57"
$ws.Range("E13").Value = "  -5.97%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.335"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.495"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008639"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").Value = "26.931.28"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.151"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "1.924.48"
$ws.Range("E24").Value = "  -8.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.843"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.058"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.125"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.27"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08841"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.949"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.428"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7201"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05259"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01923"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.406"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.929"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.171"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5169"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1630"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8588"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -15.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.179"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4814"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06252"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.76%  "
